$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Modelo NER": rename the last Precision/Recall/F1 block header from
# "Dropout_03 Classes" to "V2 Classes" and fill in the previously empty
# V/W/X (Precision/Recall/F1) columns for rows 18-46.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Modelo NER")

$ws4.Range("U17").Value = "V2 Classes"

$ws4.Range("V18").Value = 0.92857142857142805
$ws4.Range("W18").Value = 0.88930348258706404
$ws4.Range("X18").Value = 0.97146739130434701
$ws4.Range("V19").Value = 0.88715953307392903
$ws4.Range("W19").Value = 0.83823529411764697
$ws4.Range("X19").Value = 0.94214876033057804
$ws4.Range("V20").Value = 0.90640394088669896
$ws4.Range("W20").Value = 0.91089108910891003
$ws4.Range("X20").Value = 0.90196078431372495
$ws4.Range("V21").Value = 0.94
$ws4.Range("W21").Value = 0.89523809523809506
$ws4.Range("X21").Value = 0.98947368421052595
$ws4.Range("V22").Value = 0.95757575757575697
$ws4.Range("W22").Value = 0.91860465116279
$ws4.Range("X22").Value = 1
$ws4.Range("V23").Value = 1
$ws4.Range("W23").Value = 1
$ws4.Range("X23").Value = 1
$ws4.Range("V24").Value = 0.93690851735015701
$ws4.Range("W24").Value = 0.88922155688622695
$ws4.Range("X24").Value = 0.99
$ws4.Range("V25").Value = 0.92105263157894701
$ws4.Range("W25").Value = 0.875
$ws4.Range("X25").Value = 0.97222222222222199
$ws4.Range("V26").Value = 0.96644295302013405
$ws4.Range("W26").Value = 0.94117647058823495
$ws4.Range("X26").Value = 0.99310344827586206
$ws4.Range("V27").Value = 0.59824046920821095
$ws4.Range("W27").Value = 0.43404255319148899
$ws4.Range("X27").Value = 0.96226415094339601
$ws4.Range("V28").Value = 0.92385786802030401
$ws4.Range("W28").Value = 0.86666666666666603
$ws4.Range("X28").Value = 0.98913043478260798
$ws4.Range("V29").Value = 0.133333333333333
$ws4.Range("W29").Value = 0.090909090909090898
$ws4.Range("X29").Value = 0.25
$ws4.Range("V30").Value = 0.91176470588235203
$ws4.Range("W30").Value = 0.874608150470219
$ws4.Range("X30").Value = 0.95221843003412898
$ws4.Range("V31").Value = 0.94527363184079505
$ws4.Range("W31").Value = 0.94059405940593999
$ws4.Range("X31").Value = 0.95
$ws4.Range("V32").Value = 0.92901802133137101
$ws4.Range("W32").Value = 0.89321074964639302
$ws4.Range("X32").Value = 0.96781609195402296
$ws4.Range("V33").Value = 0.90146750524108998
$ws4.Range("W33").Value = 0.83657587548638102
$ws4.Range("X33").Value = 0.97727272727272696
$ws4.Range("V34").Value = 0.94533762057877802
$ws4.Range("W34").Value = 0.91874999999999996
$ws4.Range("X34").Value = 0.97350993377483397
$ws4.Range("V35").Value = 0.82456140350877105
$ws4.Range("W35").Value = 0.71573604060913698
$ws4.Range("X35").Value = 0.972413793103448
$ws4.Range("V36").Value = 0.92307692307692302
$ws4.Range("W36").Value = 0.85714285714285698
$ws4.Range("X36").Value = 1
$ws4.Range("V37").Value = 0.875
$ws4.Range("W37").Value = 0.77777777777777701
$ws4.Range("X37").Value = 1
$ws4.Range("V38").Value = 0.96321070234113704
$ws4.Range("W38").Value = 0.95049504950495001
$ws4.Range("X38").Value = 0.97627118644067801
$ws4.Range("V39").Value = 0.42857142857142799
$ws4.Range("W39").Value = 1
$ws4.Range("X39").Value = 0.27272727272727199
$ws4.Range("V40").Value = 0.93687707641196005
$ws4.Range("W40").Value = 0.91558441558441495
$ws4.Range("X40").Value = 0.95918367346938704
$ws4.Range("V41").Value = 0.61395348837209296
$ws4.Range("W41").Value = 0.44337811900191898
$ws4.Range("X41").Value = 0.99784017278617698
$ws4.Range("V42").Value = 0.88973966309341501
$ws4.Range("W42").Value = 0.80248618784530301
$ws4.Range("X42").Value = 0.99828178694157998
$ws4.Range("V43").Value = 0
$ws4.Range("W43").Value = 0
$ws4.Range("X43").Value = 0
$ws4.Range("V44").Value = 0.91130012150668205
$ws4.Range("W44").Value = 0.84383438343834305
$ws4.Range("X44").Value = 0.99049128367670303
$ws4.Range("V45").Value = 0.932668329177057
$ws4.Range("W45").Value = 0.89473684210526305
$ws4.Range("X45").Value = 0.97395833333333304
$ws4.Range("V46").Value = 0.99158007498660905
$ws4.Range("W46").Value = 0.99899845451170299
$ws4.Range("X46").Value = 0.98427105840719997

$ws4.Activate()
$ws4.Range("Z17:AD46").Select()

# ---------------------------------------------------------------------------
# Sheet "Modelo Instances clustering": selection moved to A1:G11.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Modelo Instances clustering")
$ws5.Activate()
$ws5.Range("A1:G11").Select()

# ---------------------------------------------------------------------------
# Sheet "Modelo RE": selection moved to A11:D17.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Modelo RE")
$ws6.Activate()
$ws6.Range("A11:D17").Select()

# ---------------------------------------------------------------------------
# Sheet "Pipeline": selection moved to the full A:C columns, and it is no
# longer the active/selected tab.
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Pipeline")
$ws7.Activate()
$ws7.Range("A:C").Select()

# ---------------------------------------------------------------------------
# Sheet "PetroOntoVec": a new (otherwise empty) row 44 is appended with a
# single formatted cell K44, extending the used range to A1:P44. This sheet
# ends up being the active tab/selection of the workbook.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PetroOntoVec")
$ws2.Range("K1").Copy()
$ws2.Range("K44").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Activate()
$ws2.Range("C48").Select()
